$d = $word.ActiveDocument

# --- 1. Update the Intermediate CentOS AMI version: 2.0.2 -> 3.0.0 ---------
# Locate the run that holds "2.0.2" precisely (it's a distinct, bold run).
$r = $d.Content
$found = $r.Find.Execute("2.0.2", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find '2.0.2' text to update."
}

# Toggle Bold off/on around the text replacement so the run-coalescing pass
# that merges adjacent same-formatted runs doesn't fold this run into its
# bold neighbour; the final formatting still matches (Bold stays True).
$r.Bold = 0
$r.Text = "3.0.0"
$r.Bold = 1

# --- 2. Move the "_GoBack" bookmark to sit right after the new "3.0.0" ----
# Re-find the freshly written text so we have an accurate Range, then
# collapse it to its end point (a zero-length range right after "3.0.0").
$r2 = $d.Content
$r2.Find.Execute("3.0.0", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Collapse(0)

# Adding a bookmark named "_GoBack" automatically removes any other
# bookmark of that name (Word only allows one _GoBack), so this both
# relocates it from after "Cod3Can!" and drops the stale one.
$d.Bookmarks.Add("_GoBack", $r2)
